# Scheduled-runner style refresh of market price/profit columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2793.75
$ws.Range("I64").Value = 2788.889
$ws.Range("J64").Value = 2800
$ws.Range("K64").Value = 2788.889
$ws.Range("L64").Value = 2800
$ws.Range("M64").Value = -2540.889
$ws.Range("N64").Value = -3296

$ws.Range("H67").Value = 2793.75
$ws.Range("I67").Value = 2788.889
$ws.Range("J67").Value = 2800
$ws.Range("K67").Value = 2788.889
$ws.Range("L67").Value = 2800
$ws.Range("M67").Value = -1930.889
$ws.Range("N67").Value = -4516

$ws.Range("H132").Value = 358689.47
$ws.Range("I132").Value = 251011.25
$ws.Range("J132").Value = 627885.0600000001
$ws.Range("K132").Value = 753033.75
$ws.Range("L132").Value = 1883655.18
$ws.Range("M132").Value = -750503.75
$ws.Range("N132").Value = -1888715.18

$ws.Range("H137").Value = 3727.2593
$ws.Range("I137").Value = 1771.4375
$ws.Range("J137").Value = 6572.091
$ws.Range("K137").Value = 5314.3125
$ws.Range("L137").Value = 19716.273
$ws.Range("M137").Value = -2764.3125
$ws.Range("N137").Value = -24816.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1028.275
$ws.Range("I61").Value = 803.34283
$ws.Range("J61").Value = 2602.8
$ws.Range("K61").Value = 803.34283
$ws.Range("L61").Value = 2602.8
$ws.Range("M61").Value = -591.34283
$ws.Range("N61").Value = -3026.8

$ws.Range("H74").Value = 5102.5
$ws.Range("I74").Value = 5267.0586
$ws.Range("J74").Value = 4543
$ws.Range("K74").Value = 5267.0586
$ws.Range("L74").Value = 4543
$ws.Range("M74").Value = -4393.0586
$ws.Range("N74").Value = -6291

$ws.Range("H77").Value = 5102.5
$ws.Range("I77").Value = 5267.0586
$ws.Range("J77").Value = 4543
$ws.Range("K77").Value = 26335.293
$ws.Range("L77").Value = 22715
$ws.Range("M77").Value = -21967.293
$ws.Range("N77").Value = -31451

$ws.Range("H136").Value = 1028.275
$ws.Range("I136").Value = 803.34283
$ws.Range("J136").Value = 2602.8
$ws.Range("K136").Value = 2410.02849
$ws.Range("L136").Value = 7808.400000000001
$ws.Range("M136").Value = 139.9715099999999
$ws.Range("N136").Value = -12908.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 60780
$ws.Range("J57").Value = 60780
$ws.Range("L57").Value = 60780
$ws.Range("N57").Value = -62220

$ws.Range("H58").Value = 96500
$ws.Range("J58").Value = 96500
$ws.Range("L58").Value = 96500
$ws.Range("N58").Value = -97088

$ws.Range("H60").Value = 19390
$ws.Range("J60").Value = 19390
$ws.Range("L60").Value = 19390
$ws.Range("N60").Value = -20588

$ws.Range("H132").Value = 44991.54
$ws.Range("J132").Value = 44991.54
$ws.Range("L132").Value = 44991.54
$ws.Range("N132").Value = -55111.54

$ws.Range("H133").Value = 50755
$ws.Range("J133").Value = 50755
$ws.Range("L133").Value = 50755
$ws.Range("N133").Value = -60875

$ws.Range("H135").Value = 43775
$ws.Range("J135").Value = 43775
$ws.Range("L135").Value = 43775
$ws.Range("N135").Value = -53915

$ws.Range("H136").Value = 60780
$ws.Range("J136").Value = 60780
$ws.Range("L136").Value = 60780
$ws.Range("N136").Value = -70980

$ws.Range("H138").Value = 41415.652
$ws.Range("J138").Value = 41415.652
$ws.Range("L138").Value = 41415.652
$ws.Range("N138").Value = -51695.652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13891559
$ws.Range("I31").Value = 1503.2273
$ws.Range("K31").Value = 1503.2273
$ws.Range("M31").Value = -1208.2273

$ws.Range("H34").Value = 13891559
$ws.Range("I34").Value = 1503.2273
$ws.Range("K34").Value = 1503.2273
$ws.Range("M34").Value = -1301.2273

$ws.Range("H47").Value = 49000
$ws.Range("J47").Value = 49000
$ws.Range("L47").Value = 49000
$ws.Range("N47").Value = -50132

$ws.Range("H58").Value = 1637.0706
$ws.Range("I58").Value = 1455.4722
$ws.Range("J58").Value = 2642.8462
$ws.Range("K58").Value = 1455.4722
$ws.Range("L58").Value = 2642.8462
$ws.Range("M58").Value = -1252.4722
$ws.Range("N58").Value = -3048.8462

$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880

$ws.Range("H134").Value = 1962.2285
$ws.Range("J134").Value = 3318.6875
$ws.Range("L134").Value = 9956.0625
$ws.Range("N134").Value = -15026.0625

$ws.Range("H136").Value = 1637.0706
$ws.Range("I136").Value = 1455.4722
$ws.Range("J136").Value = 2642.8462
$ws.Range("K136").Value = 4366.4166
$ws.Range("L136").Value = 7928.5386
$ws.Range("M136").Value = -1816.4166
$ws.Range("N136").Value = -13028.5386

$ws.Range("H138").Value = 45840
$ws.Range("J138").Value = 45840
$ws.Range("L138").Value = 45840
$ws.Range("N138").Value = -56120

$ws.Range("H140").Value = 80429.28999999999
$ws.Range("J140").Value = 80429.28999999999
$ws.Range("L140").Value = 80429.28999999999
$ws.Range("N140").Value = -90789.28999999999

$ws.Range("H141").Value = 34900
$ws.Range("J141").Value = 34900
$ws.Range("L141").Value = 34900
$ws.Range("N141").Value = -45260

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7043069
$ws.Range("I131").Value = 83333640
$ws.Range("J131").Value = 862.5231
$ws.Range("K131").Value = 250000920
$ws.Range("L131").Value = 2587.5693
$ws.Range("M131").Value = -249995880
$ws.Range("N131").Value = -12667.5693

$ws.Range("H137").Value = 2647.0908
$ws.Range("I137").Value = 683.3333
$ws.Range("J137").Value = 4006.6155
$ws.Range("K137").Value = 2049.9999
$ws.Range("L137").Value = 12019.8465
$ws.Range("M137").Value = 3050.0001
$ws.Range("N137").Value = -22219.8465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 5581671
$ws.Range("I14").Value = 7800140
$ws.Range("K14").Value = 7800140
$ws.Range("M14").Value = -7799972

$ws.Range("H114").Value = 40000
$ws.Range("J114").Value = 40000
$ws.Range("L114").Value = 40000
$ws.Range("N114").Value = -48678

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9978.628000000001
$ws.Range("I132").Value = 11892.956
$ws.Range("K132").Value = 35678.868
$ws.Range("M132").Value = -33148.868

$ws.Range("H136").Value = 3600
$ws.Range("I136").Value = 1815.3846
$ws.Range("K136").Value = 5446.1538
$ws.Range("M136").Value = -2896.1538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996

$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984
